$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1234
$ws.Range("E3").Value = 1234
$ws.Range("E4").Value = 5678
$ws.Range("E5").Value = 5678
$ws.Range("E6").Value = 5678
$ws.Range("E7").Value = 1234

$ws.Range("E7").Select()
